$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Make room for the new "localdb" command-group in the #system sheet.
#    - Insert a new column at N (14): this pushes the existing "macro"
#      through "xml" blocks of data one column to the right (N..AC -> O..AD).
#    - Make room in the alphabetically sorted "target" list (column A) for
#      "localdb" between "json" and "macro" by shifting A14:A29 down to
#      A15:A30 one cell at a time (keeps the shift scoped to column A only).
# ---------------------------------------------------------------------------
$ws.Columns.Item(14).Insert()

for ($r = 29; $r -ge 14; $r--) {
    $srcValue = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $srcValue
}
$ws.Range("A14").Value = "localdb"

# ---------------------------------------------------------------------------
# 2. Populate the new "localdb" column (N) with its header and commands.
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------------
# 3. Fix up the workbook-level defined names so they keep pointing at the
#    right data after the column insert / row insert above, and register
#    the new "localdb" name.
# ---------------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")

Write-Host "localdb command group added"
